# Weekly price-sheet update: insert a new weekly record as row 14,
# pushing the existing rows 14-43 down to 15-44 (dimension grows to A1:R44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14 - shifts old rows 14..43 down to 15..44
# and carries the row's formatting (incl. the date-number-format style on D).
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Cells.Item(14, 1).Value  = 1
$ws.Cells.Item(14, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value  = 45070
$ws.Cells.Item(14, 5).Value  = 15
$ws.Cells.Item(14, 6).Value  = 100112044
$ws.Cells.Item(14, 7).Value  = "Perejil"
$ws.Cells.Item(14, 8).Value  = "Sin especificar"
$ws.Cells.Item(14, 9).Value  = "Primera"
$ws.Cells.Item(14, 10).Value = 270
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1250
$ws.Cells.Item(14, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 625
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(14, 18).Value = "Hortaliza"
